# Applies the weekly refresh of the "Sandia" (watermelon) price rows.
# Each existing record's price block (rows 35..84) is shifted down by three
# rows (i.e. every block of "Extra/Primera/Segunda" rows is replaced with the
# data that used to sit three rows below it), three brand-new rows are
# inserted at the top of the shifted range (35..37, date 2022-01-20, Región
# del Maule) and the range grows by three rows overall (35..87) because the
# last three historical rows (what used to be rows 83..84, now new rows
# 86..87) are preserved instead of falling off the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function NewDate($y, $m, $d) {
    return Get-Date -Year $y -Month $m -Day $d -Hour 0 -Minute 0 -Second 0
}

# Columns that stay constant for every data row in this sheet.
$MercadoId = 11
$Mercado   = "Vega Monumental Concepción"
$Region    = "Bíobío"
$Codreg    = 8
$CatId     = 100112028
$Categoria = "Sandia"
$Variedad  = "Sin especificar"
$KgOUnidad = 1
$Clasif    = "Hortaliza"

# Full target content for rows 35..87 (the range that changes).
# Fields: Row, Year, Month, Day, Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm, Unidad, Origen, PrecioKg
$rows = @(
    @{Row=35; Y=2022; M=1;  D=20; I="Extra";   J=400;  K=2500; L=2500; M2=2500; N="$/unidad"; O="Región del Maule";   P=2500},
    @{Row=36; Y=2022; M=1;  D=20; I="Primera"; J=400;  K=2000; L=2000; M2=2000; N="$/unidad"; O="Región del Maule";   P=2000},
    @{Row=37; Y=2022; M=1;  D=20; I="Segunda"; J=400;  K=1500; L=1500; M2=1500; N="$/unidad"; O="Región del Maule";   P=1500},

    @{Row=38; Y=2021; M=2;  D=4;  I="Extra";   J=500;  K=3000; L=3000; M2=3000; N="$/unidad"; O="Región de O'Higgins"; P=3000},
    @{Row=39; Y=2021; M=2;  D=4;  I="Primera"; J=500;  K=2500; L=2500; M2=2500; N="$/unidad"; O="Región de O'Higgins"; P=2500},
    @{Row=40; Y=2021; M=2;  D=4;  I="Segunda"; J=500;  K=2000; L=2000; M2=2000; N="$/unidad"; O="Región de O'Higgins"; P=2000},

    @{Row=41; Y=2022; M=1;  D=18; I="Extra";   J=500;  K=3000; L=3000; M2=3000; N="$/unidad"; O="Región de O'Higgins"; P=3000},
    @{Row=42; Y=2022; M=1;  D=18; I="Primera"; J=800;  K=2500; L=2500; M2=2500; N="$/unidad"; O="Región de O'Higgins"; P=2500},
    @{Row=43; Y=2022; M=1;  D=18; I="Segunda"; J=800;  K=2000; L=2000; M2=2000; N="$/unidad"; O="Región de O'Higgins"; P=2000},

    @{Row=44; Y=2021; M=1;  D=6;  I="Extra";   J=500;  K=3000; L=3000; M2=3000; N="$/unidad"; O="Región de O'Higgins"; P=3000},
    @{Row=45; Y=2021; M=1;  D=6;  I="Primera"; J=500;  K=2500; L=2500; M2=2500; N="$/unidad"; O="Región de O'Higgins"; P=2500},
    @{Row=46; Y=2021; M=1;  D=6;  I="Segunda"; J=500;  K=2000; L=2000; M2=2000; N="$/unidad"; O="Región de O'Higgins"; P=2000},

    @{Row=47; Y=2020; M=12; D=22; I="Extra";   J=400;  K=3200; L=3200; M2=3200; N="$/unidad"; O="Región de O'Higgins"; P=3200},
    @{Row=48; Y=2020; M=12; D=22; I="Primera"; J=400;  K=2800; L=2800; M2=2800; N="$/unidad"; O="Región de O'Higgins"; P=2800},
    @{Row=49; Y=2020; M=12; D=22; I="Segunda"; J=400;  K=2400; L=2400; M2=2400; N="$/unidad"; O="Región de O'Higgins"; P=2400},

    @{Row=50; Y=2021; M=12; D=16; I="Primera"; J=450;  K=2500; L=3000; M2=2778; N="$/unidad"; O="Región de O'Higgins"; P=2778},

    @{Row=51; Y=2021; M=1;  D=19; I="Extra";   J=400;  K=3000; L=3000; M2=3000; N="$/unidad"; O="Región de O'Higgins"; P=3000},
    @{Row=52; Y=2021; M=1;  D=19; I="Primera"; J=400;  K=2500; L=2500; M2=2500; N="$/unidad"; O="Región de O'Higgins"; P=2500},
    @{Row=53; Y=2021; M=1;  D=19; I="Segunda"; J=300;  K=2000; L=2000; M2=2000; N="$/unidad"; O="Región de O'Higgins"; P=2000},

    @{Row=54; Y=2022; M=1;  D=12; I="Primera"; J=2700; K=2000; L=2200; M2=2089; N="$/unidad"; O="Región de O'Higgins"; P=2089},

    @{Row=55; Y=2021; M=1;  D=8;  I="Extra";   J=500;  K=3500; L=3500; M2=3500; N="$/unidad"; O="Región de O'Higgins"; P=3500},
    @{Row=56; Y=2021; M=1;  D=8;  I="Primera"; J=500;  K=3000; L=3000; M2=3000; N="$/unidad"; O="Región de O'Higgins"; P=3000},
    @{Row=57; Y=2021; M=1;  D=8;  I="Segunda"; J=500;  K=2500; L=2500; M2=2500; N="$/unidad"; O="Región de O'Higgins"; P=2500},

    @{Row=58; Y=2022; M=1;  D=13; I="Extra";   J=400;  K=3000; L=3000; M2=3000; N="$/unidad"; O="Región de O'Higgins"; P=3000},
    @{Row=59; Y=2022; M=1;  D=13; I="Primera"; J=400;  K=2500; L=2500; M2=2500; N="$/unidad"; O="Región de O'Higgins"; P=2500},
    @{Row=60; Y=2022; M=1;  D=13; I="Segunda"; J=400;  K=2000; L=2000; M2=2000; N="$/unidad"; O="Región de O'Higgins"; P=2000},

    @{Row=61; Y=2021; M=2;  D=24; I="Extra";   J=300;  K=3000; L=3000; M2=3000; N="$/unidad"; O="Región de O'Higgins"; P=3000},
    @{Row=62; Y=2021; M=2;  D=24; I="Primera"; J=500;  K=2500; L=2500; M2=2500; N="$/unidad"; O="Región de O'Higgins"; P=2500},
    @{Row=63; Y=2021; M=2;  D=24; I="Segunda"; J=500;  K=2000; L=2000; M2=2000; N="$/unidad"; O="Región de O'Higgins"; P=2000},

    @{Row=64; Y=2021; M=3;  D=17; I="Primera"; J=300;  K=2500; L=2500; M2=2500; N="$/unidad"; O="Región de O'Higgins"; P=2500},
    @{Row=65; Y=2021; M=3;  D=17; I="Segunda"; J=300;  K=2000; L=2000; M2=2000; N="$/unidad"; O="Región de O'Higgins"; P=2000},

    @{Row=66; Y=2021; M=1;  D=15; I="Extra";   J=500;  K=3500; L=3500; M2=3500; N="$/unidad"; O="Región de O'Higgins"; P=3500},
    @{Row=67; Y=2021; M=1;  D=15; I="Primera"; J=500;  K=3000; L=3000; M2=3000; N="$/unidad"; O="Región de O'Higgins"; P=3000},
    @{Row=68; Y=2021; M=1;  D=15; I="Segunda"; J=500;  K=2500; L=2500; M2=2500; N="$/unidad"; O="Región de O'Higgins"; P=2500},

    @{Row=69; Y=2020; M=12; D=31; I="Extra";   J=400;  K=3000; L=3000; M2=3000; N="$/unidad"; O="Región de O'Higgins"; P=3000},
    @{Row=70; Y=2020; M=12; D=31; I="Primera"; J=400;  K=2500; L=2500; M2=2500; N="$/unidad"; O="Región de O'Higgins"; P=2500},
    @{Row=71; Y=2020; M=12; D=31; I="Segunda"; J=400;  K=2000; L=2000; M2=2000; N="$/unidad"; O="Región de O'Higgins"; P=2000},

    @{Row=72; Y=2021; M=1;  D=12; I="Extra";   J=500;  K=3500; L=3500; M2=3500; N="$/kilo (volumen en unidades)"; O="Región de O'Higgins"; P=3500},
    @{Row=73; Y=2021; M=1;  D=12; I="Primera"; J=500;  K=3000; L=3000; M2=3000; N="$/kilo (volumen en unidades)"; O="Región de O'Higgins"; P=3000},
    @{Row=74; Y=2021; M=1;  D=12; I="Segunda"; J=500;  K=2500; L=2500; M2=2500; N="$/kilo (volumen en unidades)"; O="Región de O'Higgins"; P=2500},

    @{Row=75; Y=2021; M=12; D=15; I="Primera"; J=1100; K=2500; L=3000; M2=2773; N="$/unidad"; O="Región de O'Higgins"; P=2773},
    @{Row=76; Y=2021; M=12; D=15; I="Segunda"; J=350;  K=2300; L=2500; M2=2414; N="$/unidad"; O="Región de O'Higgins"; P=2414},

    @{Row=77; Y=2020; M=12; D=29; I="Extra";   J=400;  K=3000; L=3000; M2=3000; N="$/unidad"; O="Región de O'Higgins"; P=3000},
    @{Row=78; Y=2020; M=12; D=29; I="Primera"; J=500;  K=2500; L=2500; M2=2500; N="$/unidad"; O="Región de O'Higgins"; P=2500},
    @{Row=79; Y=2020; M=12; D=29; I="Segunda"; J=400;  K=2000; L=2000; M2=2000; N="$/unidad"; O="Región de O'Higgins"; P=2000},

    @{Row=80; Y=2021; M=2;  D=9;  I="Extra";   J=500;  K=3000; L=3000; M2=3000; N="$/unidad"; O="Región de O'Higgins"; P=3000},
    @{Row=81; Y=2021; M=2;  D=9;  I="Primera"; J=500;  K=2500; L=2500; M2=2500; N="$/unidad"; O="Región de O'Higgins"; P=2500},
    @{Row=82; Y=2021; M=2;  D=9;  I="Segunda"; J=500;  K=2000; L=2000; M2=2000; N="$/unidad"; O="Región de O'Higgins"; P=2000},

    @{Row=83; Y=2021; M=3;  D=9;  I="Extra";   J=300;  K=2800; L=2800; M2=2800; N="$/unidad"; O="Región de O'Higgins"; P=2800},
    @{Row=84; Y=2021; M=3;  D=9;  I="Primera"; J=300;  K=2500; L=2500; M2=2500; N="$/unidad"; O="Región de O'Higgins"; P=2500},
    @{Row=85; Y=2021; M=3;  D=9;  I="Segunda"; J=300;  K=2200; L=2200; M2=2200; N="$/unidad"; O="Región de O'Higgins"; P=2200},

    @{Row=86; Y=2021; M=11; D=25; I="Primera"; J=200;  K=700;  L=800;  M2=750;  N="$/kilo (volumen en unidades)"; O="Perú"; P=750},
    @{Row=87; Y=2021; M=11; D=11; I="Primera"; J=600;  K=800;  L=900;  M2=850;  N="$/kilo (volumen en unidades)"; O="Perú"; P=850}
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value  = $MercadoId
    $ws.Cells.Item($row, 2).Value  = $Mercado
    $ws.Cells.Item($row, 3).Value  = $Region
    $ws.Cells.Item($row, 4).Value  = (NewDate $r.Y $r.M $r.D)
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value  = $Codreg
    $ws.Cells.Item($row, 6).Value  = $CatId
    $ws.Cells.Item($row, 7).Value  = $Categoria
    $ws.Cells.Item($row, 8).Value  = $Variedad
    $ws.Cells.Item($row, 9).Value  = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M2
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $KgOUnidad
    $ws.Cells.Item($row, 18).Value = $Clasif
}
